$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B95").Value = 79244
$ws.Range("B96").Value = 79244
$ws.Range("B97").Value = 79244
$ws.Range("AX97").Value = "Anders Esplund, Enviro Planning, Anna Sjövall, Pia Edfors, Sofia Berg"
$ws.Range("AX97").ClearFormats()
$ws.Range("B98").Value = 98931
$ws.Range("B99").Value = 98931
$ws.Range("A100").Value = 130964543
$ws.Range("B100").Value = 57881
$ws.Range("E100").Value = 100049
$ws.Range("F100").Value = "Spillkråka"
$ws.Range("F100").ClearFormats()
$ws.Range("G100").Value = "Dryocopus martius"
$ws.Range("G100").ClearFormats()
$ws.Range("H100").Value = "(Linnaeus, 1758)"
$ws.Range("H100").ClearFormats()
$ws.Range("Q100").Value = 509622
$ws.Range("R100").Value = 6718933
$ws.Range("AC100").Value = "Födosökspår . inventering åt vasa vind"
$ws.Range("AC100").ClearFormats()
$ws.Range("A101").Value = 130964531
$ws.Range("B101").Value = 79244
$ws.Range("E101").Value = 6425
$ws.Range("F101").Value = "Garnlav"
$ws.Range("F101").ClearFormats()
$ws.Range("G101").Value = "Alectoria sarmentosa"
$ws.Range("G101").ClearFormats()
$ws.Range("H101").Value = "(Ach.) Ach."
$ws.Range("H101").ClearFormats()
$ws.Range("Q101").Value = 509889
$ws.Range("R101").Value = 6719134
$ws.Range("AC101").Value = "Flera . inventering åt vasa vind"
$ws.Range("AC101").ClearFormats()
$ws.Range("B102").Value = 79244
$ws.Range("A104").Value = 130964545
$ws.Range("B104").Value = 57073
$ws.Range("D104").Value = "LC"
$ws.Range("D104").ClearFormats()
$ws.Range("E104").Value = 100138
$ws.Range("F104").Value = "Tjäder"
$ws.Range("F104").ClearFormats()
$ws.Range("G104").Value = "Tetrao urogallus"
$ws.Range("G104").ClearFormats()
$ws.Range("H104").Value = "Linnaeus, 1758"
$ws.Range("H104").ClearFormats()
$ws.Range("Q104").Value = 509535
$ws.Range("R104").Value = 6718925
$ws.Range("AC104").Value = "Spillning . inventering åt vasa vind"
$ws.Range("AC104").ClearFormats()
$ws.Range("A105").Value = 130964541
$ws.Range("B105").Value = 91809
$ws.Range("D105").Value = "NT"
$ws.Range("D105").ClearFormats()
$ws.Range("E105").Value = 1202
$ws.Range("F105").Value = "Ullticka"
$ws.Range("F105").ClearFormats()
$ws.Range("G105").Value = "Phellinidium ferrugineofuscum"
$ws.Range("G105").ClearFormats()
$ws.Range("H105").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("H105").ClearFormats()
$ws.Range("Q105").Value = 509703
$ws.Range("R105").Value = 6719018
$ws.Range("AC105").Value = "Enstaka . inventering åt vasa vind"
$ws.Range("AC105").ClearFormats()
$ws.Range("A106").Value = 130964537
$ws.Range("B106").Value = 79244
$ws.Range("D106").Value = "NT"
$ws.Range("D106").ClearFormats()
$ws.Range("E106").Value = 6425
$ws.Range("F106").Value = "Garnlav"
$ws.Range("F106").ClearFormats()
$ws.Range("G106").Value = "Alectoria sarmentosa"
$ws.Range("G106").ClearFormats()
$ws.Range("H106").Value = "(Ach.) Ach."
$ws.Range("H106").ClearFormats()
$ws.Range("Q106").Value = 509822
$ws.Range("R106").Value = 6718960
$ws.Range("AC106").Value = "Rikligt . inventering åt vasa vind"
$ws.Range("AC106").ClearFormats()
$ws.Range("A107").Value = 130964642
$ws.Range("B107").Value = 99037
$ws.Range("D107").Value = "LC"
$ws.Range("D107").ClearFormats()
$ws.Range("E107").Value = 221952
$ws.Range("F107").Value = "Spindelblomster"
$ws.Range("F107").ClearFormats()
$ws.Range("G107").Value = "Neottia cordata"
$ws.Range("G107").ClearFormats()
$ws.Range("H107").Value = "(L.) Rich."
$ws.Range("H107").ClearFormats()
$ws.Range("Q107").Value = 509917
$ws.Range("R107").Value = 6719042
$ws.Range("AC107").Value = "Måttliga förekomster . inventering åt vasa vind"
$ws.Range("AC107").ClearFormats()
$ws.Range("B108").Value = 98931
$ws.Range("A109").Value = 130964574
$ws.Range("B109").Value = 79244
$ws.Range("E109").Value = 6425
$ws.Range("F109").Value = "Garnlav"
$ws.Range("F109").ClearFormats()
$ws.Range("G109").Value = "Alectoria sarmentosa"
$ws.Range("G109").ClearFormats()
$ws.Range("H109").Value = "(Ach.) Ach."
$ws.Range("H109").ClearFormats()
$ws.Range("Q109").Value = 509667
$ws.Range("R109").Value = 6719184
$ws.Range("A110").Value = 130964647
$ws.Range("B110").Value = 92107
$ws.Range("E110").Value = 658
$ws.Range("F110").Value = "Rosenticka"
$ws.Range("F110").ClearFormats()
$ws.Range("G110").Value = "Fomitopsis rosea"
$ws.Range("G110").ClearFormats()
$ws.Range("H110").Value = "(Alb. & Schwein.:Fr.) P.Karst."
$ws.Range("H110").ClearFormats()
$ws.Range("Q110").Value = 509741
$ws.Range("R110").Value = 6718998
$ws.Range("A112").Value = 130964533
$ws.Range("B112").Value = 79244
$ws.Range("D112").Value = "NT"
$ws.Range("D112").ClearFormats()
$ws.Range("E112").Value = 6425
$ws.Range("F112").Value = "Garnlav"
$ws.Range("F112").ClearFormats()
$ws.Range("G112").Value = "Alectoria sarmentosa"
$ws.Range("G112").ClearFormats()
$ws.Range("H112").Value = "(Ach.) Ach."
$ws.Range("H112").ClearFormats()
$ws.Range("Q112").Value = 509984
$ws.Range("R112").Value = 6719028
$ws.Range("AC112").Value = "Rikligt . inventering åt vasa vind"
$ws.Range("AC112").ClearFormats()
$ws.Range("A113").Value = 130964650
$ws.Range("B113").Value = 92268
$ws.Range("D113").Value = "VU"
$ws.Range("D113").ClearFormats()
$ws.Range("E113").Value = 1209
$ws.Range("F113").Value = "Rynkskinn"
$ws.Range("F113").ClearFormats()
$ws.Range("G113").Value = "Hermanssonia centrifuga"
$ws.Range("G113").ClearFormats()
$ws.Range("H113").Value = "(P. Karst.) Zmitr."
$ws.Range("H113").ClearFormats()
$ws.Range("Q113").Value = 509694
$ws.Range("R113").Value = 6718936
$ws.Range("A114").Value = 130964645
$ws.Range("B114").Value = 99037
$ws.Range("D114").Value = "LC"
$ws.Range("D114").ClearFormats()
$ws.Range("E114").Value = 221952
$ws.Range("F114").Value = "Spindelblomster"
$ws.Range("F114").ClearFormats()
$ws.Range("G114").Value = "Neottia cordata"
$ws.Range("G114").ClearFormats()
$ws.Range("H114").Value = "(L.) Rich."
$ws.Range("H114").ClearFormats()
$ws.Range("Q114").Value = 509804
$ws.Range("R114").Value = 6719024
$ws.Range("AC114").Value = "Måttliga förekomster . inventering åt vasa vind"
$ws.Range("AC114").ClearFormats()
$ws.Range("B115").Value = 79244
$ws.Range("B116").Value = 98931
$ws.Range("B117").Value = 79244
$ws.Range("B118").Value = 98931
$ws.Range("B119").Value = 92268
$ws.Range("B120").Value = 98931
$ws.Range("B121").Value = 79244
$ws.Range("A122").Value = 130964640
$ws.Range("B122").Value = 57881
$ws.Range("D122").Value = "NT"
$ws.Range("D122").ClearFormats()
$ws.Range("E122").Value = 100049
$ws.Range("F122").Value = "Spillkråka"
$ws.Range("F122").ClearFormats()
$ws.Range("G122").Value = "Dryocopus martius"
$ws.Range("G122").ClearFormats()
$ws.Range("H122").Value = "(Linnaeus, 1758)"
$ws.Range("H122").ClearFormats()
$ws.Range("Q122").Value = 509697
$ws.Range("R122").Value = 6719144
$ws.Range("AC122").Value = "Gamla födosöksspår . inventering åt vasa vind"
$ws.Range("AC122").ClearFormats()
$ws.Range("A123").Value = 130964542
$ws.Range("B123").Value = 57073
$ws.Range("D123").Value = "LC"
$ws.Range("D123").ClearFormats()
$ws.Range("E123").Value = 100138
$ws.Range("F123").Value = "Tjäder"
$ws.Range("F123").ClearFormats()
$ws.Range("G123").Value = "Tetrao urogallus"
$ws.Range("G123").ClearFormats()
$ws.Range("H123").Value = "Linnaeus, 1758"
$ws.Range("H123").ClearFormats()
$ws.Range("Q123").Value = 509635
$ws.Range("R123").Value = 6718941
$ws.Range("AC123").Value = "Spillning . inventering åt vasa vind"
$ws.Range("AC123").ClearFormats()
$ws.Range("A124").Value = 130964644
$ws.Range("B124").Value = 98918
$ws.Range("E124").Value = 220093
$ws.Range("F124").Value = "Korallrot"
$ws.Range("F124").ClearFormats()
$ws.Range("G124").Value = "Corallorhiza trifida"
$ws.Range("G124").ClearFormats()
$ws.Range("H124").Value = "Châtel."
$ws.Range("H124").ClearFormats()
$ws.Range("Q124").Value = 509801
$ws.Range("R124").Value = 6719017
$ws.Range("AC124").Value = "Sparsamma förekomster . inventering åt vasa vind"
$ws.Range("AC124").ClearFormats()
$ws.Range("A125").Value = 130964390
$ws.Range("B125").Value = 99014
$ws.Range("D125").Value = "VU"
$ws.Range("D125").ClearFormats()
$ws.Range("E125").Value = 220787
$ws.Range("F125").Value = "Knärot"
$ws.Range("F125").ClearFormats()
$ws.Range("G125").Value = "Goodyera repens"
$ws.Range("G125").ClearFormats()
$ws.Range("H125").Value = "(L.) R. Br."
$ws.Range("H125").ClearFormats()
$ws.Range("Q125").Value = 509475
$ws.Range("R125").Value = 6718881
$ws.Range("Y125").NumberFormat = "@"
$ws.Range("Y125").Value = "2025-07-03"
$ws.Range("Y125").ClearFormats()
$ws.Range("AA125").NumberFormat = "@"
$ws.Range("AA125").Value = "2025-07-03"
$ws.Range("AA125").ClearFormats()
$ws.Range("AC125").Value = "Måttliga förekomster, Ca 10-15 plantor . inventering åt vasa vind"
$ws.Range("AC125").ClearFormats()
$ws.Range("A126").Value = 130964643
$ws.Range("B126").Value = 98931
$ws.Range("D126").Value = "LC"
$ws.Range("D126").ClearFormats()
$ws.Range("E126").Value = 219790
$ws.Range("F126").Value = "Fläcknycklar"
$ws.Range("F126").ClearFormats()
$ws.Range("G126").Value = "Dactylorhiza maculata"
$ws.Range("G126").ClearFormats()
$ws.Range("H126").Value = "(L.) Soó"
$ws.Range("H126").ClearFormats()
$ws.Range("Q126").Value = 509829
$ws.Range("R126").Value = 6719000
$ws.Range("Y126").NumberFormat = "@"
$ws.Range("Y126").Value = "2025-07-02"
$ws.Range("Y126").ClearFormats()
$ws.Range("AA126").NumberFormat = "@"
$ws.Range("AA126").Value = "2025-07-02"
$ws.Range("AA126").ClearFormats()
$ws.Range("AC126").Value = "Måttlig förekomst . inventering åt vasa vind"
$ws.Range("AC126").ClearFormats()
$ws.Range("B127").Value = 92504
$ws.Range("B128").Value = 79244
$ws.Range("A129").Value = 130964646
$ws.Range("B129").Value = 91809
$ws.Range("D129").Value = "NT"
$ws.Range("D129").ClearFormats()
$ws.Range("E129").Value = 1202
$ws.Range("F129").Value = "Ullticka"
$ws.Range("F129").ClearFormats()
$ws.Range("G129").Value = "Phellinidium ferrugineofuscum"
$ws.Range("G129").ClearFormats()
$ws.Range("H129").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("H129").ClearFormats()
$ws.Range("Q129").Value = 509764
$ws.Range("R129").Value = 6719043
$ws.Range("AC129").Value = "Måttliga förekomster . inventering åt vasa vind"
$ws.Range("AC129").ClearFormats()
$ws.Range("A130").Value = 130964639
$ws.Range("B130").Value = 57073
$ws.Range("D130").Value = "LC"
$ws.Range("D130").ClearFormats()
$ws.Range("E130").Value = 100138
$ws.Range("F130").Value = "Tjäder"
$ws.Range("F130").ClearFormats()
$ws.Range("G130").Value = "Tetrao urogallus"
$ws.Range("G130").ClearFormats()
$ws.Range("H130").Value = "Linnaeus, 1758"
$ws.Range("H130").ClearFormats()
$ws.Range("Q130").Value = 509645
$ws.Range("R130").Value = 6719169
$ws.Range("AC130").Value = "Vinterspillning . inventering åt vasa vind"
$ws.Range("AC130").ClearFormats()
$ws.Range("B131").Value = 79244
